# Update cryptocurrency price/volume figures (scrape refresh, GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.329.04'
$ws.Range('E2').Value = '  +2.15%  '
$ws.Range('D3').Value = '1.662.67'
$ws.Range('E3').Value = '  +1.30%  '
$ws.Range('E4').Value = '  -0.34%  '
$ws.Range('D5').Value = "'220.31"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.12%  '
$ws.Range('D6').Value = "'0.508"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.76%  '
$ws.Range('E7').Value = '  -0.39%  '
$ws.Range('E8').Value = '  +1.28%  '
$ws.Range('E9').Value = '  +0.29%  '
$ws.Range('D10').Value = "'19.97"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.26%  '
$ws.Range('E11').Value = '  +0.97%  '
$ws.Range('D12').Value = '1.894.10'
$ws.Range('E12').Value = '  +1.24%  '
$ws.Range('D13').Value = '1.673.34'
$ws.Range('E13').Value = '  +1.86%  '
$ws.Range('E14').Value = '  +1.19%  '
$ws.Range('E15').Value = '  +1.74%  '
$ws.Range('E16').Value = '  +4.19%  '
$ws.Range('D17').Value = '27.310.94'
$ws.Range('E17').Value = '  +2.09%  '
$ws.Range('E18').Value = '  +0.37%  '
$ws.Range('D19').Value = "'223.06"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.55%  '
$ws.Range('D21').Value = "'6.80"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +9.09%  '
$ws.Range('E22').Value = '  +1.79%  '
$ws.Range('D23').Value = "'2.49"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.70%  '
$ws.Range('D24').Value = "'9.29"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.09%  '
$ws.Range('D25').Value = "'147.49"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.25%  '
$ws.Range('E26').Value = '  -0.29%  '
$ws.Range('E27').Value = '  +3.94%  '
$ws.Range('E28').Value = '  +0.88%  '
$ws.Range('D29').Value = "'16.09"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.97%  '
$ws.Range('D30').Value = "'0.0515"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.49%  '
$ws.Range('E31').Value = '  +0.73%  '
$ws.Range('E32').Value = '  +0.71%  '
$ws.Range('E33').Value = '  +0.33%  '
$ws.Range('E34').Value = '  +2.26%  '
$ws.Range('D35').Value = '1.262.85'
$ws.Range('E35').Value = '  -1.94%  '
$ws.Range('E36').Value = '  +0.13%  '
$ws.Range('E37').Value = '  -0.30%  '
$ws.Range('E38').Value = '  +0.29%  '
$ws.Range('E39').Value = '  +1.94%  '
$ws.Range('E40').Value = '  -0.31%  '
$ws.Range('D41').Value = "'0.815"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.24%  '
$ws.Range('E42').Value = '  +2.03%  '
$ws.Range('D43').Value = '1.805.97'
$ws.Range('E43').Value = '  +1.39%  '
$ws.Range('E44').Value = '  -4.14%  '
$ws.Range('D45').Value = "'61.91"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.91%  '
$ws.Range('D46').Value = "'92.45"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.63%  '
$ws.Range('E47').Value = '  +1.44%  '
$ws.Range('D48').Value = "'0.0517"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.04%  '
$ws.Range('D49').Value = "'0.0986"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.01%  '
$ws.Range('E51').Value = '  +0.47%  '
